$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 536 all hold the "Förändrad" date serial value.
# They currently contain 45204 (2023-10-05) and must become 45205 (2023-10-06).
$ws.Range("C2:C536").Value = 45205
